$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A79").Value = "Web"
$ws.Range("B79").Value = "https://www.trinitycollege.com/local-trinity/UK/drama-speech-dance/exam-help"

$ws.Range("G48").Select()
